# Updated cryptos list data: refresh Price (D) and Volume(1h) (E) columns
# for rows 2-51, plus a coin-order swap for rows 34/35 (Filecoin <-> ARBITRUM).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 34 and 35 first: coin identity (name/link) swap with refreshed price data
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.159"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.54%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.506"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.23%  "

# Remaining rows: Price (D) and Volume(1h) (E) refresh only
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.407.10"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.84"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.49"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4663"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3727"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07401"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8908"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07961"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.19"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.884.92"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.434"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.620"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.72"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008961"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +4.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.94"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.436.64"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.164"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.61"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.074.16"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.56"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.863"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.60"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.096"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.171"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.49"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08916"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7555"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.971"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.600"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.086"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05302"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01958"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.991"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.169"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5218"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1647"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.364"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4914"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.38"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.56"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.643"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06280"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.90"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.97%  "
